$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67.8255122017956
$ws.Range("L2").Value = 73.9047667329172

$ws.Range("B3").Value = 55.3352948193206
$ws.Range("L3").Value = 48.6313846276599

$ws.Range("B4").Value = 40.563158268296
$ws.Range("L4").Value = 36.9618773883661

$ws.Range("B6").Value = 64.7210522905015

$ws.Range("B7").Value = 67.1354010141054
$ws.Range("C7").Value = 73.1501666372061
$ws.Range("L7").Value = 63.442019211072

$ws.Range("B8").Value = 68.4899159160604
$ws.Range("L8").Value = 60.1821932205212

$ws.Range("B9").Value = 61.8798692282585
$ws.Range("L9").Value = 54.0183622108344
